$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New branch "mod-schedule" row: add rows 4-7 to the schedule sheet.

# Row 4
$ws.Range("A4").Value = "sadas"
$ws.Range("B4").Value = "sa"
$ws.Range("F4").Value = 123

# Row 5
$ws.Range("A5").Value = "Jim"
$ws.Range("B5").Value = "Jim"
$ws.Range("H5").Value = 123123
$ws.Range("K5").Value = "Opportunity"
$ws.Range("L5").Value = "Won"

# Row 6
$ws.Range("A6").Value = "Loki"
$ws.Range("B6").Value = "Lyon"
$ws.Range("F6").Value = 3123
$ws.Range("G6").Value = 123
$ws.Range("H6").Value = 123
$ws.Range("K6").Value = "Opportunity"
$ws.Range("L6").Value = "Won"

# Row 7
$ws.Range("A7").Value = "NEW "
$ws.Range("B7").Value = "oLd"
$ws.Range("F7").Value = "asdas"
$ws.Range("G7").Value = 24352435
$ws.Range("K7").Value = "Opportunity"
$ws.Range("L7").Value = "Won"

# Move the active selection to F7 to match the final workbook state.
$ws.Range("F7").Select()
